# Add three new columns (D, E, F) of header metadata to Sheet1, matching
# the "ORG_GRAW_IDENOLD" / "ORG_GRAW_IDENNEW" / "ORG_GRAW_STATUS" headers
# that were added to the workbook upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (these also grow the shared-string table and the
# worksheet dimension/row spans automatically).
$ws.Range("D1").Value = "ORG_GRAW_IDENOLD"
$ws.Range("E1").Value = "ORG_GRAW_IDENNEW"
$ws.Range("F1").Value = "ORG_GRAW_STATUS"

# Match the left-aligned header style already used by A1:C1.
$ws.Range("D1:F1").HorizontalAlignment = -4131

# Give the new columns the same wide data-entry width as the rest of the
# table (about 28 characters wide).
$ws.Range("D1:F22").ColumnWidth = 27.166666666666668

# Leave the cursor/selection on the newly added column, like the author's
# Excel session did before saving.
$null = $ws.Range("F3").Select()
